$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 167, shifting rows 167-176 down to 168-177.
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new data record.
$ws.Cells.Item(167, 1).Value = 9
$ws.Cells.Item(167, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(167, 3).Value = "Metropolitana"
$ws.Cells.Item(167, 4).Value = 44610
$ws.Cells.Item(167, 5).Value = 13
$ws.Cells.Item(167, 6).Value = "Fruta"
$ws.Cells.Item(167, 7).Value = 100101
$ws.Cells.Item(167, 8).Value = "Berries"
$ws.Cells.Item(167, 9).Value = 100101001
$ws.Cells.Item(167, 10).Value = "Arándano (blue)"
$ws.Cells.Item(167, 11).Value = "Sin especificar"
$ws.Cells.Item(167, 12).Value = "Primera"
$ws.Cells.Item(167, 13).Value = 280
$ws.Cells.Item(167, 14).Value = 3600
$ws.Cells.Item(167, 15).Value = 3600
$ws.Cells.Item(167, 16).Value = 3600
$ws.Cells.Item(167, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(167, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(167, 19).Value = 1800
$ws.Cells.Item(167, 20).Value = 2
